$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.899.06'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.877.80'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D4").Value = '''1.018'
$ws.Range("E4").Value = '  +1.48%  '
$ws.Range("D5").Value = '''334.68'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("D7").Value = '''0.4696'
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").Value = '''46.86'
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("D10").Value = '''0.07949'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").Value = '''1.007'
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("D12").Value = '''21.60'
$ws.Range("D13").Value = '1.876.98'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = '''5.949'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '''7.102'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '''0.06789'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '''17.03'
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("E21").Value = '  +1.31%  '
$ws.Range("D22").Value = '27.898.50'
$ws.Range("D23").Value = '''5.467'
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").Value = '''2.359'
$ws.Range("E25").Value = '  +2.57%  '
$ws.Range("D26").Value = '2.099.37'
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").Value = '''159.53'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").Value = '''19.89'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").Value = '''2.078'
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").Value = '''120.80'
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("D32").Value = '''0.09529'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").Value = '''0.9579'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '''3.655'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("D35").Value = '''5.314'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  -7.23%  '
$ws.Range("D37").Value = '''0.06102'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '''0.02243'
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").Value = '''1.203'
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").Value = '''1.015'
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").Value = '''8.107'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '''0.5902'
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").Value = '''0.1891'
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '''10.20'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").Value = '''1.270'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").Value = '''0.5647'
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("D47").Value = '''12.14'
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("D48").Value = '''3.391'
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Value = '''1.917'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").Value = '''113.46'
$ws.Range("E51").Value = '  +1.05%  '

# Reset style on text-forced price cells so no quotePrefix style sticks (keeps original default styling)
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
